$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the new I (I0) and J (IF) columns, styled to match the
# existing bold/bordered/centered header style used by columns B through H.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data values for columns I (I0) and J (IF), row by row.
$data = @{
    2 = @(7, 7)
    3 = @(5, 5)
    4 = @(7, 7)
    5 = @(6, 6)
    6 = @(3, 3)
    7 = @(7, 7)
    8 = @(8, 8)
    9 = @(8, 8)
    10 = @(7, 7)
    11 = @(6, 7)
    12 = @(7, 7)
    13 = @(5, 5)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(7, 7)
    17 = @(4, 4)
    18 = @(8, 8)
    19 = @(8, 8)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(6, 6)
    23 = @(6, 6)
    24 = @(7, 8)
    25 = @(7, 7)
    26 = @(7, 7)
    27 = @(5, 6)
    28 = @(7, 7)
    29 = @(7, 7)
    30 = @(4, 5)
    31 = @(7, 7)
    32 = @(7, 8)
    33 = @(7, 8)
    34 = @(8, 9)
    35 = @(6, 7)
    36 = @(7, 7)
    37 = @(7, 8)
    38 = @(3, 3)
    39 = @(5, 5)
    40 = @(7, 7)
    41 = @(6, 7)
    42 = @(7, 7)
    43 = @(8, 8)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
